$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(1)
$f = $sh.TextFrame.TextRange.Font
Write-Host "Size:" $f.Size
$f.Size = $f.Size
